$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Solar (column E) capacity values for "Open year" 2022 and 2024
# to reflect updated data from upstream processes through 2024.
$ws.Range("E24").Value = 241.114
$ws.Range("E26").Value = 216.705
